$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Cells.Item(2, 2).Value = 'Bitcoin'
$ws.Cells.Item(2, 3).Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Cells.Item(2, 4).Value = '65.889.22'
$ws.Cells.Item(2, 5).Value = '  +0.32%  '

$ws.Cells.Item(3, 2).Value = 'Ethereum'
$ws.Cells.Item(3, 3).Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Cells.Item(3, 4).Value = '2.664.61'
$ws.Cells.Item(3, 5).Value = '  -0.37%  '

$ws.Cells.Item(4, 2).Value = 'TetherUSD'
$ws.Cells.Item(4, 3).Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
Set-TextValue $ws.Cells.Item(4, 4) '1.00'
$ws.Cells.Item(4, 5).Value = '  -0.06%  '

$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
Set-TextValue $ws.Cells.Item(5, 4) '600.30'
$ws.Cells.Item(5, 5).Value = '  +0.21%  '

$ws.Cells.Item(6, 2).Value = 'Solana'
$ws.Cells.Item(6, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue $ws.Cells.Item(6, 4) '160.83'
$ws.Cells.Item(6, 5).Value = '  +2.78%  '

$ws.Cells.Item(7, 2).Value = 'XRP'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
Set-TextValue $ws.Cells.Item(7, 4) '0.643'
$ws.Cells.Item(7, 5).Value = '  +4.05%  '

$ws.Cells.Item(8, 2).Value = 'USDC'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextValue $ws.Cells.Item(8, 4) '1.00'
$ws.Cells.Item(8, 5).Value = '  -0.05%  '

$ws.Cells.Item(9, 2).Value = 'Dogecoin'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextValue $ws.Cells.Item(9, 4) '0.127'
$ws.Cells.Item(9, 5).Value = '  -1.93%  '

$ws.Cells.Item(10, 2).Value = 'Cardano'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Cells.Item(10, 4) '0.402'
$ws.Cells.Item(10, 5).Value = '  +0.46%  '

$ws.Cells.Item(11, 2).Value = 'Toncoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextValue $ws.Cells.Item(11, 4) '5.89'
$ws.Cells.Item(11, 5).Value = '  +0.42%  '

$ws.Cells.Item(12, 2).Value = 'TRON'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Cells.Item(12, 4) '0.157'
$ws.Cells.Item(12, 5).Value = '  +1.63%  '

$ws.Cells.Item(13, 2).Value = 'Avalanche'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue $ws.Cells.Item(13, 4) '29.27'
$ws.Cells.Item(13, 5).Value = '  -0.08%  '

$ws.Cells.Item(14, 2).Value = 'ShibaInu'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Cells.Item(14, 4) '0.0000196'
$ws.Cells.Item(14, 5).Value = '  -0.41%  '

$ws.Cells.Item(15, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(15, 4).Value = '3.143.34'
$ws.Cells.Item(15, 5).Value = '  -0.34%  '

$ws.Cells.Item(16, 2).Value = 'WrappedBTC'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(16, 4).Value = '65.761.15'
$ws.Cells.Item(16, 5).Value = '  +0.33%  '

$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '2.650.97'
$ws.Cells.Item(17, 5).Value = '  -1.75%  '

$ws.Cells.Item(18, 2).Value = 'Chainlink'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Cells.Item(18, 4) '12.63'
$ws.Cells.Item(18, 5).Value = '  -1.99%  '

$ws.Cells.Item(19, 2).Value = 'Polkadot'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Cells.Item(19, 4) '4.82'
$ws.Cells.Item(19, 5).Value = '  +0.62%  '

$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue $ws.Cells.Item(20, 4) '356.98'
$ws.Cells.Item(20, 5).Value = '  +1.47%  '

$ws.Cells.Item(21, 2).Value = 'Uniswap'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue $ws.Cells.Item(21, 4) '7.49'
$ws.Cells.Item(21, 5).Value = '  -0.72%  '

$ws.Cells.Item(22, 2).Value = 'Dai'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws.Cells.Item(22, 4) '0.999'
$ws.Cells.Item(22, 5).Value = '  -0.06%  '

$ws.Cells.Item(23, 2).Value = 'Litecoin'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue $ws.Cells.Item(23, 4) '70.05'
$ws.Cells.Item(23, 5).Value = '  +0.34%  '

$ws.Cells.Item(24, 2).Value = 'SuiNetwork'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
Set-TextValue $ws.Cells.Item(24, 4) '1.81'
$ws.Cells.Item(24, 5).Value = '  +10.27%  '

$ws.Cells.Item(25, 2).Value = 'PEPE'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue $ws.Cells.Item(25, 4) '0.0000114'
$ws.Cells.Item(25, 5).Value = '  +2.80%  '

$ws.Cells.Item(26, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Cells.Item(26, 4) '9.76'
$ws.Cells.Item(26, 5).Value = '  +1.69%  '

$ws.Cells.Item(27, 2).Value = 'Fetch.AI'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Cells.Item(27, 4) '1.63'
$ws.Cells.Item(27, 5).Value = '  +2.57%  '

$ws.Cells.Item(28, 2).Value = 'Bittensor'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Cells.Item(28, 4) '581.85'
$ws.Cells.Item(28, 5).Value = '  +10.67%  '

$ws.Cells.Item(29, 2).Value = 'Aptos'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Cells.Item(29, 4) '8.18'
$ws.Cells.Item(29, 5).Value = '  +2.02%  '

$ws.Cells.Item(30, 2).Value = 'Kaspa'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Cells.Item(30, 4) '0.164'
$ws.Cells.Item(30, 5).Value = '  -1.22%  '

$ws.Cells.Item(31, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Cells.Item(31, 4) '0.999'
$ws.Cells.Item(31, 5).Value = '  -0.20%  '

$ws.Cells.Item(32, 2).Value = 'PancakeSwap'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue $ws.Cells.Item(32, 4) '2.16'
$ws.Cells.Item(32, 5).Value = '  +1.18%  '

$ws.Cells.Item(33, 2).Value = 'ImmutableX'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Cells.Item(33, 4) '1.83'
$ws.Cells.Item(33, 5).Value = '  +3.59%  '

$ws.Cells.Item(34, 2).Value = 'RenderToken'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Cells.Item(34, 4) '6.76'
$ws.Cells.Item(34, 5).Value = '  +4.87%  '

$ws.Cells.Item(35, 2).Value = 'NEARProtocol'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Cells.Item(35, 4) '5.51'
$ws.Cells.Item(35, 5).Value = '  +1.10%  '

$ws.Cells.Item(36, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Cells.Item(36, 4) '0.423'
$ws.Cells.Item(36, 5).Value = '  +0.16%  '

$ws.Cells.Item(37, 2).Value = 'EthereumClassic'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Cells.Item(37, 4) '20.67'
$ws.Cells.Item(37, 5).Value = '  +0.34%  '

$ws.Cells.Item(38, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Cells.Item(38, 4) '1.00'
$ws.Cells.Item(38, 5).Value = '  -0.04%  '

$ws.Cells.Item(39, 2).Value = 'Stacks'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Cells.Item(39, 4) '1.97'
$ws.Cells.Item(39, 5).Value = '  +1.98%  '

$ws.Cells.Item(40, 2).Value = 'Monero'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Cells.Item(40, 4) '154.11'
$ws.Cells.Item(40, 5).Value = '  -2.46%  '

$ws.Cells.Item(41, 2).Value = 'dogwifhat'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Cells.Item(41, 4) '2.55'
$ws.Cells.Item(41, 5).Value = '  +11.26%  '

$ws.Cells.Item(42, 2).Value = 'Aave'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Cells.Item(42, 4) '162.45'
$ws.Cells.Item(42, 5).Value = '  -0.84%  '

$ws.Cells.Item(43, 2).Value = 'Filecoin'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Cells.Item(43, 4) '4.12'
$ws.Cells.Item(43, 5).Value = '  +0.08%  '

$ws.Cells.Item(44, 2).Value = 'Hedera'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Cells.Item(44, 4) '0.0620'
$ws.Cells.Item(44, 5).Value = '  +1.85%  '

$ws.Cells.Item(45, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Cells.Item(45, 4) '23.57'
$ws.Cells.Item(45, 5).Value = '  +3.61%  '

$ws.Cells.Item(46, 2).Value = 'Mantle'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Cells.Item(46, 4) '0.646'
$ws.Cells.Item(46, 5).Value = '  +0.91%  '

$ws.Cells.Item(47, 2).Value = 'VeChain'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Cells.Item(47, 4) '0.0260'
$ws.Cells.Item(47, 5).Value = '  +1.09%  '

$ws.Cells.Item(48, 2).Value = 'Stellar'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Cells.Item(48, 4) '0.102'
$ws.Cells.Item(48, 5).Value = '  +1.76%  '

$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Cells.Item(49, 4) '19.84'
$ws.Cells.Item(49, 5).Value = '  -1.42%  '

$ws.Cells.Item(50, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(50, 4).Value = '0.0₆0247'
$ws.Cells.Item(50, 5).Value = '  -6.34%  '

$ws.Cells.Item(51, 2).Value = 'ONDO'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
Set-TextValue $ws.Cells.Item(51, 4) '0.820'
$ws.Cells.Item(51, 5).Value = '  +1.11%  '

Write-Output "done"